# Apply "Complete public law data collection and optimize code" update
# to the public_law_word_count_by_congre worksheet.
#
# Congress.gov re-pull refreshed the Page Count / Word Count / Public Law
# Count figures for most Congresses and finished out the still-partial
# 118th Congress row; the stray in-progress 119th Congress row (which had
# only a handful of laws recorded) is removed now that it is no longer
# part of the published series.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Refreshed counts for the 94th-103rd Congresses ---
$ws.Range("B2").Value = 4119
$ws.Range("C2").Value = 1959666
$ws.Range("D2").Value = 588

$ws.Range("B3").Value = 5403
$ws.Range("C3").Value = 2615715
$ws.Range("D3").Value = 633

$ws.Range("B4").Value = 4948
$ws.Range("C4").Value = 2331477
$ws.Range("D4").Value = 613

$ws.Range("B5").Value = 4339
$ws.Range("C5").Value = 1982634
$ws.Range("D5").Value = 473

$ws.Range("B6").Value = 4889
$ws.Range("C6").Value = 2187411
$ws.Range("D6").Value = 623

$ws.Range("B7").Value = 7205
$ws.Range("C7").Value = 3282060
$ws.Range("D7").Value = 666

$ws.Range("B8").Value = 7851
$ws.Range("C8").Value = 3513194
$ws.Range("D8").Value = 713

$ws.Range("B9").Value = 8287
$ws.Range("C9").Value = 3670798
$ws.Range("D9").Value = 650

$ws.Range("B10").Value = 7540
$ws.Range("C10").Value = 3249657
$ws.Range("D10").Value = 590

$ws.Range("B11").Value = 7547
$ws.Range("C11").Value = 3248009
$ws.Range("D11").Value = 465

# --- Refreshed counts for the 113th Congress ---
$ws.Range("B21").Value = 5342
$ws.Range("C21").Value = 2388640
$ws.Range("D21").Value = 296

# --- Refreshed counts for the 116th-118th Congresses ---
$ws.Range("B24").Value = 8501
$ws.Range("C24").Value = 3827126
$ws.Range("D24").Value = 344

$ws.Range("B25").Value = 8809
$ws.Range("C25").Value = 3947924
$ws.Range("D25").Value = 362

$ws.Range("B26").Value = 4428
$ws.Range("C26").Value = 1941440
$ws.Range("D26").Value = 274

# --- Drop the stray, still-incomplete 119th Congress row ---
$ws.Rows.Item(27).Delete()

# --- Match the author's final selection state ---
$ws.Range("C26").Select()
